# Apply the edit described by the diff:
#  1. Header row (A1:K1) becomes bold.
#  2. The "target" column (G2:G13) values change from "deuteron" to "d".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the header row bold.
$ws.Range("A1:K1").Font.Bold = $true

# Update the "target" column values from "deuteron" to "d" for all data rows.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = "d"
}
